# Add a new "Levorg" learning-record row (row 4) to the 学習記録 sheet and
# bring the whole used range's alignment in line with the rest of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlTop = -4160

# --- Row 3: the date cell becomes free text ("12/19/") but keeps its date
#     number format (mirrors what Excel does when you overtype a
#     date-formatted cell with plain text). ---
$ws.Range("A3").Value = "12/19/"

# --- New row 4: Levorg PC entry added 12/20 ---
$ws.Range("A4").NumberFormat = "mm/dd/yy"
$ws.Range("A4").Value = "12/20/2023"

$ws.Range("B4").Value = "LevorgPC"

$ws.Range("C4").Value = "PPO ヒンジ報酬，ドアノブ報酬のみ"
$ws.Range("C4").Characters(1, 4).Font.Name = "Arial"
$ws.Range("C4").Characters(1, 4).Font.Size = 10
$ws.Range("C4").Characters(5, 14).Font.Name = "Noto Sans CJK SC"
$ws.Range("C4").Characters(5, 14).Font.Size = 10

$ws.Range("D4").Value = "報酬のスケールを小さくした"

$ws.Range("E4").Value = "ハンドルの付け根を引っ張って開けている様子"

$F4Text = "なんか惜しいので，ドアヒンジfrictionとdampng上げる，ドアノブ報酬上げる，ドアノブ付け根urdf修正`nロボット自体のdampingと速度limit考えたほうが良いかもしれない`n1000epochくらいでこの動きが生まれ始めた様子"
$ws.Range("F4").Value = $F4Text
$ws.Range("F4").WrapText = $true

# Rich-text runs inside F4 (Arial for the Latin/technical terms, Noto Sans
# CJK SC for the Japanese prose) -- offsets computed against $F4Text.
$ws.Range("F4").Characters(1, 14).Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").Characters(1, 14).Font.Size = 10
$ws.Range("F4").Characters(15, 8).Font.Name = "Arial"
$ws.Range("F4").Characters(15, 8).Font.Size = 10
$ws.Range("F4").Characters(23, 1).Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").Characters(23, 1).Font.Size = 10
$ws.Range("F4").Characters(24, 6).Font.Name = "Arial"
$ws.Range("F4").Characters(24, 6).Font.Size = 10
$ws.Range("F4").Characters(30, 21).Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").Characters(30, 21).Font.Size = 10
$ws.Range("F4").Characters(51, 4).Font.Name = "Arial"
$ws.Range("F4").Characters(51, 4).Font.Size = 10
$ws.Range("F4").Characters(55, 10).Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").Characters(55, 10).Font.Size = 10
$ws.Range("F4").Characters(65, 7).Font.Name = "Arial"
$ws.Range("F4").Characters(65, 7).Font.Size = 10
$ws.Range("F4").Characters(72, 3).Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").Characters(72, 3).Font.Size = 10
$ws.Range("F4").Characters(75, 5).Font.Name = "Arial"
$ws.Range("F4").Characters(75, 5).Font.Size = 10
$ws.Range("F4").Characters(80, 15).Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").Characters(80, 15).Font.Size = 10
$ws.Range("F4").Characters(95, 9).Font.Name = "Arial"
$ws.Range("F4").Characters(95, 9).Font.Size = 10
$ws.Range("F4").Characters(104, 17).Font.Name = "Noto Sans CJK SC"
$ws.Range("F4").Characters(104, 17).Font.Size = 10

$ws.Range("G4").Value = "重み保存済み"

# --- Alignment: the whole used range (including the new row) is now
#     centre/top instead of general/bottom. ---
$ws.Range("A1:G4").HorizontalAlignment = $xlCenter
$ws.Range("A1:G4").VerticalAlignment = $xlTop

# --- Row heights ---
$ws.Rows(2).RowHeight = 15.45
$ws.Rows(3).RowHeight = 15.45
$ws.Rows(4).RowHeight = 33.55

# --- Column widths (character units; engine persists width+5/6) ---
$ws.Columns("A:B").ColumnWidth = 10.74
$ws.Columns("C:C").ColumnWidth = 33.62
$ws.Columns("D:D").ColumnWidth = 37.8
$ws.Columns("E:E").ColumnWidth = 39.26
$ws.Columns("F:F").ColumnWidth = 93.77
$ws.Columns("G:G").ColumnWidth = 18.82

# --- View: selection moved to E13, first visible column is B ---
$ws.Range("E13").Select()
